$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data table (rows 2..289) is being refreshed: a brand-new pair of
# rows (the latest market day) is inserted right after the header by
# overwriting rows 185-186 with new data, and every row from 187 onward is
# pushed down by 2 rows. The two rows that fall off the bottom (old 288-289)
# are appended as new rows 290-291.

# Make sure the new rows inherit the same date number format as the rest of
# column D before the bulk copy populates them (so no new style gets minted).
$ws.Range("D290:D291").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Shift rows 187..289 down to 189..291 (read first so the write doesn't
# clobber data still being read).
$block = $ws.Range("A187:R289").Value()
$ws.Range("A189:R291").Value = $block

# Row 185: new "Morada(o)" record for the latest market date.
$ws.Cells.Item(185, 1).Value = 4
$ws.Cells.Item(185, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(185, 3).Value = "Los Lagos"
$ws.Cells.Item(185, 4).Value = 44460
$ws.Cells.Item(185, 5).Value = 10
$ws.Cells.Item(185, 6).Value = 100112004
$ws.Cells.Item(185, 7).Value = "Cebolla"
$ws.Cells.Item(185, 8).Value = "Morada(o)"
$ws.Cells.Item(185, 9).Value = "1a (guarda)"
$ws.Cells.Item(185, 10).Value = 100
$ws.Cells.Item(185, 11).Value = 13000
$ws.Cells.Item(185, 12).Value = 13000
$ws.Cells.Item(185, 13).Value = 13000
$ws.Cells.Item(185, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(185, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(185, 16).Value = 722
$ws.Cells.Item(185, 17).Value = 18
$ws.Cells.Item(185, 18).Value = "Hortaliza"

# Row 186: new "Sin especificar" record for the latest market date.
$ws.Cells.Item(186, 1).Value = 4
$ws.Cells.Item(186, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value = "Los Lagos"
$ws.Cells.Item(186, 4).Value = 44460
$ws.Cells.Item(186, 5).Value = 10
$ws.Cells.Item(186, 6).Value = 100112004
$ws.Cells.Item(186, 7).Value = "Cebolla"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "1a (guarda)"
$ws.Cells.Item(186, 10).Value = 300
$ws.Cells.Item(186, 11).Value = 7000
$ws.Cells.Item(186, 12).Value = 7000
$ws.Cells.Item(186, 13).Value = 7000
$ws.Cells.Item(186, 14).Value = "`$/malla 16 kilos"
$ws.Cells.Item(186, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(186, 16).Value = 438
$ws.Cells.Item(186, 17).Value = 16
$ws.Cells.Item(186, 18).Value = "Hortaliza"
